$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# B6 previously held the text "proceso" (shared string 16). It now holds "Leido"
# (the same shared string used elsewhere, e.g. B2/B3/B5).
$ws.Range("B6").Value = "Leido"

# The old "proceso" text is renamed to "Preceso" and moved down to the new cell B7
# (row for "6-Desarrollo del Prototipo").
$ws.Range("B7").Value = "Preceso"

# Update the active selection to the new cell B7.
$ws.Activate()
$ws.Range("B7").Select()
